$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new file path entry in the next empty row (row 10)
$ws.Range("A10").Value = "D:/My Files/My Personal Stuff/Documents"

# Move the active selection to B14 as recorded in the saved workbook view
$ws.Range("B14").Select()
